# Auto-generated Excel COM-interop script applying scheduled-runner value updates
# to the Kraken_Profits "Leve profit" worksheets (one table per job class).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 19.5
$ws.Range("I6").Value = 19.5
$ws.Range("K6").Value = 58.5
$ws.Range("M6").Value = 53.5
$ws.Range("H52").Value = 249.5
$ws.Range("I52").Value = 249.5
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 748.5
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -588.5
$ws.Range("N52").ClearContents()
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 4250
$ws.Range("K74").Value = 4250
$ws.Range("M74").Value = -3314
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 4250
$ws.Range("K77").Value = 21250
$ws.Range("M77").Value = -16570
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 2500
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 7500
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -4950
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 3397.875
$ws.Range("I138").Value = 2768.3
$ws.Range("K138").Value = 8304.900000000001
$ws.Range("M138").Value = -3164.900000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2225.5715
$ws.Range("I2").Value = 1072
$ws.Range("K2").Value = 1072
$ws.Range("M2").Value = -959
$ws.Range("H61").Value = 5895
$ws.Range("I61").Value = 5371.25
$ws.Range("K61").Value = 5371.25
$ws.Range("M61").Value = -5159.25
$ws.Range("H74").Value = 4794
$ws.Range("I74").Value = 2733.4119
$ws.Range("K74").Value = 2733.4119
$ws.Range("M74").Value = -1859.4119
$ws.Range("H77").Value = 4794
$ws.Range("I77").Value = 2733.4119
$ws.Range("K77").Value = 13667.0595
$ws.Range("M77").Value = -9299.059499999999
$ws.Range("H110").Value = 842.6667
$ws.Range("I110").Value = 842.6667
$ws.Range("K110").Value = 842.6667
$ws.Range("M110").Value = 1202.3333
$ws.Range("H116").Value = 2225.5715
$ws.Range("I116").Value = 1072
$ws.Range("K116").Value = 1072
$ws.Range("M116").Value = 1222
$ws.Range("H122").Value = 4328.4375
$ws.Range("I122").Value = 4861.727
$ws.Range("J122").Value = 3155.2
$ws.Range("K122").Value = 14585.181
$ws.Range("L122").Value = 9465.599999999999
$ws.Range("M122").Value = -12135.181
$ws.Range("N122").Value = -14365.6
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 4946.5
$ws.Range("I132").Value = 3994.3333
$ws.Range("K132").Value = 11982.9999
$ws.Range("M132").Value = -9452.999899999999
$ws.Range("H136").Value = 5895
$ws.Range("I136").Value = 5371.25
$ws.Range("K136").Value = 16113.75
$ws.Range("M136").Value = -13563.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2225.5715
$ws.Range("I3").Value = 1072
$ws.Range("K3").Value = 1072
$ws.Range("M3").Value = -958
$ws.Range("H20").Value = 1036.8572
$ws.Range("I20").Value = 769.3333
$ws.Range("J20").Value = 1237.5
$ws.Range("K20").Value = 769.3333
$ws.Range("L20").Value = 1237.5
$ws.Range("M20").Value = -522.3333
$ws.Range("N20").Value = -1731.5
$ws.Range("H134").Value = 5858.9
$ws.Range("J134").Value = 7465.6665
$ws.Range("L134").Value = 22396.9995
$ws.Range("N134").Value = -27466.9995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 193
$ws.Range("I7").Value = 124.5
$ws.Range("J7").Value = 330
$ws.Range("K7").Value = 124.5
$ws.Range("L7").Value = 330
$ws.Range("M7").Value = -11.5
$ws.Range("N7").Value = -556
$ws.Range("H31").Value = 10462.667
$ws.Range("I31").Value = 15465
$ws.Range("K31").Value = 15465
$ws.Range("M31").Value = -15170
$ws.Range("H34").Value = 10462.667
$ws.Range("I34").Value = 15465
$ws.Range("K34").Value = 15465
$ws.Range("M34").Value = -15263
$ws.Range("H95").Value = 31666.666
$ws.Range("J95").Value = 31666.666
$ws.Range("L95").Value = 31666.666
$ws.Range("N95").Value = -37158.666
$ws.Range("H105").Value = 1112.125
$ws.Range("I105").Value = 1056.7142
$ws.Range("K105").Value = 1056.7142
$ws.Range("M105").Value = 690.2858000000001
$ws.Range("H107").Value = 723
$ws.Range("J107").Value = 796
$ws.Range("L107").Value = 796
$ws.Range("N107").Value = -4636
$ws.Range("H122").Value = 1494.25
$ws.Range("I122").Value = 1461
$ws.Range("J122").Value = 1527.5
$ws.Range("K122").Value = 4383
$ws.Range("L122").Value = 4582.5
$ws.Range("M122").Value = -1933
$ws.Range("N122").Value = -9482.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 812.375
$ws.Range("I7").Value = 83.5
$ws.Range("J7").Value = 2999
$ws.Range("K7").Value = 250.5
$ws.Range("L7").Value = 8997
$ws.Range("M7").Value = -138.5
$ws.Range("N7").Value = -9221
$ws.Range("H14").Value = 632.1667
$ws.Range("I14").Value = 632.1667
$ws.Range("K14").Value = 1896.5001
$ws.Range("M14").Value = -1723.5001
$ws.Range("H81").Value = 15
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 15
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 45
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -2291
$ws.Range("H84").Value = 15
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 15
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 135
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -11367
$ws.Range("H98").Value = 55
$ws.Range("I98").Value = 55
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 165
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 1333
$ws.Range("N98").ClearContents()
$ws.Range("H114").Value = 1287.3334
$ws.Range("I114").Value = 943.3333
$ws.Range("J114").Value = 1631.3334
$ws.Range("K114").Value = 2829.9999
$ws.Range("L114").Value = 4894.0002
$ws.Range("M114").Value = 424.0001000000002
$ws.Range("N114").Value = -11402.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H97").Value = 1670.7142
$ws.Range("I97").Value = 592.75
$ws.Range("K97").Value = 592.75
$ws.Range("M97").Value = -96.75
$ws.Range("H98").Value = 4466.6665
$ws.Range("J98").Value = 4466.6665
$ws.Range("L98").Value = 4466.6665
$ws.Range("N98").Value = -10456.6665
$ws.Range("H102").Value = 1416.5
$ws.Range("I102").Value = 1416.5
$ws.Range("K102").Value = 1416.5
$ws.Range("M102").Value = 205.5
$ws.Range("H122").Value = 1266.3334
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 4025
$ws.Range("I126").Value = 4033.3333
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 12099.9999
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -9629.999899999999
$ws.Range("N126").Value = -16940

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3176.5
$ws.Range("I40").Value = 3111.8
$ws.Range("K40").Value = 3111.8
$ws.Range("M40").Value = -2975.8
$ws.Range("H46").Value = 4334.8335
$ws.Range("I46").Value = 2250
$ws.Range("J46").Value = 4751.8
$ws.Range("K46").Value = 2250
$ws.Range("L46").Value = 4751.8
$ws.Range("M46").Value = -2062
$ws.Range("N46").Value = -5127.8
$ws.Range("H55").Value = 1146.8334
$ws.Range("I55").Value = 973
$ws.Range("K55").Value = 973
$ws.Range("M55").Value = -800
$ws.Range("H68").Value = 2173.1428
$ws.Range("I68").Value = 2092.4
$ws.Range("J68").Value = 2375
$ws.Range("K68").Value = 2092.4
$ws.Range("L68").Value = 2375
$ws.Range("M68").Value = -1343.4
$ws.Range("N68").Value = -3873
$ws.Range("H71").Value = 2173.1428
$ws.Range("I71").Value = 2092.4
$ws.Range("J71").Value = 2375
$ws.Range("K71").Value = 10462
$ws.Range("L71").Value = 11875
$ws.Range("M71").Value = -6718
$ws.Range("N71").Value = -19363
$ws.Range("H132").Value = 5500
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 31625.154
$ws.Range("I2").Value = 36452.637
$ws.Range("J2").Value = 5074
$ws.Range("K2").Value = 36452.637
$ws.Range("L2").Value = 5074
$ws.Range("M2").Value = -36340.637
$ws.Range("N2").Value = -5298
$ws.Range("H107").Value = 2867.375
$ws.Range("J107").Value = 7000
$ws.Range("L107").Value = 21000
$ws.Range("N107").Value = -24840
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3765.75
$ws.Range("I122").Value = 3523.875
$ws.Range("J122").Value = 4249.5
$ws.Range("K122").Value = 10571.625
$ws.Range("L122").Value = 12748.5
$ws.Range("M122").Value = -8121.625
$ws.Range("N122").Value = -17648.5
$ws.Range("H126").Value = 1196.3334
$ws.Range("I126").Value = 1196.3334
$ws.Range("K126").Value = 3589.0002
$ws.Range("M126").Value = -1119.0002
$ws.Range("H132").Value = 7537.25
$ws.Range("I132").Value = 6499.75
$ws.Range("K132").Value = 19499.25
$ws.Range("M132").Value = -16969.25

